$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1969
$ws.Range("J17").Value = 2259.7
$ws.Range("L17").Value = 6779.099999999999
$ws.Range("N17").Value = -7115.099999999999

$ws.Range("H40").Value = 2360.3635
$ws.Range("J40").Value = 3640
$ws.Range("L40").Value = 3640
$ws.Range("N40").Value = -3990

$ws.Range("H129").Value = 1621.7894
$ws.Range("J129").Value = 2467
$ws.Range("L129").Value = 7401
$ws.Range("N129").Value = -17401

$ws.Range("H132").Value = 2145.4666
$ws.Range("I132").Value = 1603.5264
$ws.Range("K132").Value = 4810.5792
$ws.Range("M132").Value = -2280.5792

$ws.Range("H138").Value = 3106.2812
$ws.Range("I138").Value = 2380.7693
$ws.Range("J138").Value = 3291.2156
$ws.Range("K138").Value = 7142.3079
$ws.Range("L138").Value = 9873.6468
$ws.Range("M138").Value = -2002.3079
$ws.Range("N138").Value = -20153.6468

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1954.3846
$ws.Range("I45").Value = 1743.8182
$ws.Range("K45").Value = 1743.8182
$ws.Range("M45").Value = -1366.8182

$ws.Range("H61").Value = 5092.353
$ws.Range("I61").Value = 5092.353
$ws.Range("K61").Value = 5092.353
$ws.Range("M61").Value = -4880.353

$ws.Range("H88").Value = 2531.0908
$ws.Range("I88").Value = 1841.8572
$ws.Range("K88").Value = 1841.8572
$ws.Range("M88").Value = -1435.8572

$ws.Range("H91").Value = 2531.0908
$ws.Range("I91").Value = 1841.8572
$ws.Range("K91").Value = 1841.8572
$ws.Range("M91").Value = -437.8571999999999

$ws.Range("H102").Value = 3918.9375
$ws.Range("I102").Value = 2142.125
$ws.Range("J102").Value = 9249.375
$ws.Range("K102").Value = 2142.125
$ws.Range("L102").Value = 9249.375
$ws.Range("M102").Value = -520.125
$ws.Range("N102").Value = -12493.375

$ws.Range("H110").Value = 2056.625
$ws.Range("I110").Value = 868
$ws.Range("K110").Value = 868
$ws.Range("M110").Value = 1177

$ws.Range("H136").Value = 5092.353
$ws.Range("I136").Value = 5092.353
$ws.Range("K136").Value = 15277.059
$ws.Range("M136").Value = -12727.059

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 529227
$ws.Range("I86").Value = 1477.5769
$ws.Range("J86").Value = 1901375.4
$ws.Range("K86").Value = 1477.5769
$ws.Range("L86").Value = 1901375.4
$ws.Range("M86").Value = -354.5769
$ws.Range("N86").Value = -1903621.4

$ws.Range("H89").Value = 529227
$ws.Range("I89").Value = 1477.5769
$ws.Range("J89").Value = 1901375.4
$ws.Range("K89").Value = 7387.8845
$ws.Range("L89").Value = 9506877
$ws.Range("M89").Value = -1771.8845
$ws.Range("N89").Value = -9518109

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 244.15
$ws.Range("I7").Value = 215.55556
$ws.Range("K7").Value = 215.55556
$ws.Range("M7").Value = -102.55556

$ws.Range("H31").Value = 3495.7742
$ws.Range("I31").Value = 1541.2106
$ws.Range("K31").Value = 1541.2106
$ws.Range("M31").Value = -1246.2106

$ws.Range("H34").Value = 3495.7742
$ws.Range("I34").Value = 1541.2106
$ws.Range("K34").Value = 1541.2106
$ws.Range("M34").Value = -1339.2106

$ws.Range("H58").Value = 2406.6316
$ws.Range("I58").Value = 2567.7
$ws.Range("K58").Value = 2567.7
$ws.Range("M58").Value = -2364.7

$ws.Range("H99").Value = 2914.5
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").ClearContents()

$ws.Range("H126").Value = 2914.5
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()

$ws.Range("H132").Value = 2146.5557
$ws.Range("I132").Value = 2146.5557
$ws.Range("K132").Value = 6439.6671
$ws.Range("M132").Value = -3909.6671

$ws.Range("H134").Value = 2270.4707
$ws.Range("I134").Value = 2256.1875
$ws.Range("K134").Value = 6768.5625
$ws.Range("M134").Value = -4233.5625

$ws.Range("H136").Value = 2406.6316
$ws.Range("I136").Value = 2567.7
$ws.Range("K136").Value = 7703.099999999999
$ws.Range("M136").Value = -5153.099999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1456.08
$ws.Range("J113").Value = 1726.8422
$ws.Range("L113").Value = 5180.5266
$ws.Range("N113").Value = -9520.526600000001

$ws.Range("H131").Value = 1817.7333
$ws.Range("I131").Value = 1166.6666
$ws.Range("J131").Value = 1980.5
$ws.Range("K131").Value = 3499.9998
$ws.Range("L131").Value = 5941.5
$ws.Range("M131").Value = 1540.0002
$ws.Range("N131").Value = -16021.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H68").Value = 8969
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 8969
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

$ws.Range("H113").Value = 7743.7617
$ws.Range("I113").Value = 4596
$ws.Range("K113").Value = 4596
$ws.Range("M113").Value = -2426

$ws.Range("H123").Value = 25999
$ws.Range("J123").Value = 25999
$ws.Range("L123").Value = 25999
$ws.Range("N123").Value = -30899

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3917.8215
$ws.Range("I40").Value = 2872.2222
$ws.Range("K40").Value = 2872.2222
$ws.Range("M40").Value = -2736.2222

$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("M56").ClearContents()

$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()

$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

$ws.Range("H96").Value = 42000
$ws.Range("J96").Value = 42000
$ws.Range("L96").Value = 42000
$ws.Range("N96").Value = -47492

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4616.3335
$ws.Range("I81").Value = 4575.25
$ws.Range("K81").Value = 9150.5
$ws.Range("M81").Value = -8089.5

$ws.Range("H84").Value = 4616.3335
$ws.Range("I84").Value = 4575.25
$ws.Range("K84").Value = 45752.5
$ws.Range("M84").Value = -40448.5

$ws.Range("H100").Value = 1703.05
$ws.Range("I100").Value = 1518.3529
$ws.Range("K100").Value = 3036.7058
$ws.Range("M100").Value = -2495.7058

$ws.Range("H132").Value = 3540.2
$ws.Range("I132").Value = 3280.4583
$ws.Range("J132").Value = 4579.1665
$ws.Range("K132").Value = 9841.374899999999
$ws.Range("L132").Value = 13737.4995
$ws.Range("M132").Value = -7311.374899999999
$ws.Range("N132").Value = -18797.4995
